# Insert a new data row before row 12. Excel's native Insert() behaviour
# shifts all existing rows (12..115) down by one (to 13..116), duplicates
# the formatting of the row above for the newly inserted row, and updates
# the sheet's used range/dimension automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new price-quote record.
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(12, 3).Value = "La Araucanía"
$ws.Cells.Item(12, 4).Value = 44490
$ws.Cells.Item(12, 5).Value = 9
$ws.Cells.Item(12, 6).Value = 100114007
$ws.Cells.Item(12, 7).Value = "Jengibre"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 25000
$ws.Cells.Item(12, 12).Value = 25000
$ws.Cells.Item(12, 13).Value = 25000
$ws.Cells.Item(12, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(12, 15).Value = "Perú"
$ws.Cells.Item(12, 16).Value = 1923
$ws.Cells.Item(12, 17).Value = 13
$ws.Cells.Item(12, 18).Value = "Hortaliza"
